# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.167.21"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").Value = "1.603.96"
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.87%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +2.91%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("D12").Value = "1.826.49"
$ws.Range("E12").Value = "  +3.51%  "
$ws.Range("D13").Value = "1.596.89"
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").Value = "26.156.52"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").Value = "0.0₃0722"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "204.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.30%  "
$ws.Range("E21").Value = "  +3.55%  "
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("E23").Value = "  +2.93%  "
$ws.Range("E24").Value = "  +10.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("E32").Value = "  +3.66%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("E36").Value = "  +10.19%  "
$ws.Range("D37").Value = "1.117.35"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.783"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("D43").Value = "1.738.69"
$ws.Range("E43").Value = "  +3.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "0.0₇0924"
$ws.Range("E51").Value = "  -17.30%  "
